# Update cryptocurrency price (column D) and volume change (column E) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '37.840.27'
$ws.Range("E2").Value = '  -0.05%  '
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '2.084.50'
$ws.Range("E3").Value = '  -0.27%  '
$ws.Range("E4").Value = '  -0.02%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '233.53'
$ws.Range("E5").Value = '  +0.31%  '
$ws.Range("E6").Value = '  -0.13%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '59.06'
$ws.Range("E7").Value = '  +2.90%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("E9").Value = '  +1.93%  '
$ws.Range("E10").Value = '  +0.88%  '
$ws.Range("E11").Value = '  +2.71%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '2.391.87'
$ws.Range("E12").Value = '  +0.11%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '14.79'
$ws.Range("E13").Value = '  +2.28%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '21.32'
$ws.Range("E14").Value = '  +0.93%  '
$ws.Range("E15").Value = '  +1.28%  '
$ws.Range("E16").Value = '  +1.47%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '2.085.64'
$ws.Range("E17").Value = '  -0.32%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '37.773.61'
$ws.Range("E18").Value = '  -0.07%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '6.17'
$ws.Range("E19").Value = '  +0.49%  '
$ws.Range("E20").Value = '  +1.31%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '0.0₃0849'
$ws.Range("E21").Value = '  +3.24%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '228.26'
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("E23").Value = '  -0.06%  '
$ws.Range("E24").Value = '  -0.33%  '
$ws.Range("E25").Value = '  +0.42%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '9.60'
$ws.Range("E26").Value = '  +7.42%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '171.08'
$ws.Range("E27").Value = '  +0.29%  '
$ws.Range("E28").Value = '  -1.99%  '
$ws.Range("E29").Value = '  -0.95%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '19.53'
$ws.Range("E30").Value = '  +0.20%  '
$ws.Range("E31").Value = '  +2.02%  '
$ws.Range("E32").Value = '  +2.24%  '
$ws.Range("E33").Value = '  +1.34%  '
$ws.Range("E34").Value = '  +1.76%  '
$ws.Range("E35").Value = '  -0.18%  '
$ws.Range("E36").Value = '  +0.71%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '1.83'
$ws.Range("E37").Value = '  -0.09%  '
$ws.Range("E38").Value = '  -0.17%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '5.43'
$ws.Range("E39").Value = '  -0.15%  '
$ws.Range("E40").Value = '  -1.57%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '99.31'
$ws.Range("E41").Value = '  +2.03%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '17.22'
$ws.Range("E42").Value = '  +9.65%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.0218'
$ws.Range("E43").Value = '  +2.34%  '
$ws.Range("E44").Value = '  -1.29%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '1.450.03'
$ws.Range("E45").Value = '  +0.00%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '1.16'
$ws.Range("E46").Value = '  -1.00%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '4.16'
$ws.Range("E47").Value = '  +3.06%  '
$ws.Range("E48").Value = '  +1.07%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '7.37'
$ws.Range("E49").Value = '  -0.31%  '
$ws.Range("E50").Value = '  -1.21%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '2.277.16'
$ws.Range("E51").Value = '  -0.30%  '
